# Auto-generated Excel COM-interop script
# Applies the automated BRVM data refresh to "Recommandations" and "Top_YTD" sheets.

$wb = $excel.ActiveWorkbook

$wsReco = $wb.Worksheets.Item("Recommandations")

$recoData = New-Object "object[,]" 51,7
$recoData[0,0] = "SUCRIVOIRE"
$recoData[0,1] = 0
$recoData[0,2] = 4
$recoData[0,3] = 3875
$recoData[0,4] = 980
$recoData[0,5] = "🟡 Observer"
$recoData[0,6] = "➖ Neutre"
$recoData[1,0] = "BRVM - SERVICES PUBLICS"
$recoData[1,1] = 0
$recoData[1,2] = 8
$recoData[1,3] = 3371.78
$recoData[1,4] = 112.08
$recoData[1,5] = "🟡 Observer"
$recoData[1,6] = "➖ Neutre"
$recoData[2,0] = "SAFCA CI"
$recoData[2,1] = 0
$recoData[2,2] = 4
$recoData[2,3] = 2735
$recoData[2,4] = 695
$recoData[2,5] = "🟡 Observer"
$recoData[2,6] = "➖ Neutre"
$recoData[3,0] = "CFAO MOTORS CI"
$recoData[3,1] = 0
$recoData[3,2] = 4
$recoData[3,3] = 2730
$recoData[3,4] = 685
$recoData[3,5] = "🟡 Observer"
$recoData[3,6] = "➖ Neutre"
$recoData[4,0] = "BRVM - AUTRES SECTEURS"
$recoData[4,1] = 0
$recoData[4,2] = 4
$recoData[4,3] = 2642.03
$recoData[4,4] = 666.32
$recoData[4,5] = "🟡 Observer"
$recoData[4,6] = "➖ Neutre"
$recoData[5,0] = "NEI-CEDA CI"
$recoData[5,1] = 0
$recoData[5,2] = 4
$recoData[5,3] = 2365
$recoData[5,4] = 600
$recoData[5,5] = "🟡 Observer"
$recoData[5,6] = "➖ Neutre"
$recoData[6,0] = "UNIWAX CI"
$recoData[6,1] = 0
$recoData[6,2] = 4
$recoData[6,3] = 2340
$recoData[6,4] = 585
$recoData[6,5] = "🟡 Observer"
$recoData[6,6] = "➖ Neutre"
$recoData[7,0] = "SETAO CI"
$recoData[7,1] = 0
$recoData[7,2] = 4
$recoData[7,3] = 2265
$recoData[7,4] = 580
$recoData[7,5] = "🟡 Observer"
$recoData[7,6] = "➖ Neutre"
$recoData[8,0] = "AIR LIQUIDE CI"
$recoData[8,1] = 0
$recoData[8,2] = 4
$recoData[8,3] = 2155
$recoData[8,4] = 540
$recoData[8,5] = "🟡 Observer"
$recoData[8,6] = "➖ Neutre"
$recoData[9,0] = "BRVM - DISTRIBUTION"
$recoData[9,1] = 0
$recoData[9,2] = 4
$recoData[9,3] = 1498.47
$recoData[9,4] = 374.85
$recoData[9,5] = "🟡 Observer"
$recoData[9,6] = "➖ Neutre"
$recoData[10,0] = "BRVM - TRANSPORT"
$recoData[10,1] = 0
$recoData[10,2] = 4
$recoData[10,3] = 1396.44
$recoData[10,4] = 350.03
$recoData[10,5] = "🟡 Observer"
$recoData[10,6] = "➖ Neutre"
$recoData[11,0] = "BRVM - AGRICULTURE"
$recoData[11,1] = 0
$recoData[11,2] = 4
$recoData[11,3] = 1329.84
$recoData[11,4] = 333.07
$recoData[11,5] = "🟡 Observer"
$recoData[11,6] = "➖ Neutre"
$recoData[12,0] = "BRVM - INDUSTRIE"
$recoData[12,1] = 0
$recoData[12,2] = 4
$recoData[12,3] = 792.57
$recoData[12,4] = 203.97
$recoData[12,5] = "🟡 Observer"
$recoData[12,6] = "➖ Neutre"
$recoData[13,0] = "BRVM-PRINCIPAL"
$recoData[13,1] = 0
$recoData[13,2] = 4
$recoData[13,3] = 705.97
$recoData[13,4] = 177.85
$recoData[13,5] = "🟡 Observer"
$recoData[13,6] = "➖ Neutre"
$recoData[14,0] = "BRVM - CONSOMMATION DE BASE"
$recoData[14,1] = 0
$recoData[14,2] = 4
$recoData[14,3] = 701.26
$recoData[14,4] = 178.9
$recoData[14,5] = "🟡 Observer"
$recoData[14,6] = "➖ Neutre"
$recoData[15,0] = "BRVM-PRESTIGE"
$recoData[15,1] = 0
$recoData[15,2] = 4
$recoData[15,3] = 525.46
$recoData[15,4] = 132.12
$recoData[15,5] = "🟡 Observer"
$recoData[15,6] = "➖ Neutre"
$recoData[16,0] = "BRVM - INDUSTRIELS"
$recoData[16,1] = 0
$recoData[16,2] = 4
$recoData[16,3] = 519.26
$recoData[16,4] = 133.24
$recoData[16,5] = "🟡 Observer"
$recoData[16,6] = "➖ Neutre"
$recoData[17,0] = "BRVM - FINANCES"
$recoData[17,1] = 0
$recoData[17,2] = 4
$recoData[17,3] = 489.53
$recoData[17,4] = 122.75
$recoData[17,5] = "🟡 Observer"
$recoData[17,6] = "➖ Neutre"
$recoData[18,0] = "BRVM - SERVICES FINANCIERS"
$recoData[18,1] = 0
$recoData[18,2] = 4
$recoData[18,3] = 481.1
$recoData[18,4] = 120.64
$recoData[18,5] = "🟡 Observer"
$recoData[18,6] = "➖ Neutre"
$recoData[19,0] = "BRVM - ENERGIE"
$recoData[19,1] = 0
$recoData[19,2] = 4
$recoData[19,3] = 444.83
$recoData[19,4] = 111.84
$recoData[19,5] = "🟡 Observer"
$recoData[19,6] = "➖ Neutre"
$recoData[20,0] = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$recoData[20,1] = 0
$recoData[20,2] = 4
$recoData[20,3] = 425.26
$recoData[20,4] = 106.19
$recoData[20,5] = "🟡 Observer"
$recoData[20,6] = "➖ Neutre"
$recoData[21,0] = "BRVM - TELECOMMUNICATIONS"
$recoData[21,1] = 0
$recoData[21,2] = 4
$recoData[21,3] = 380.24
$recoData[21,4] = 95.34
$recoData[21,5] = "🟡 Observer"
$recoData[21,6] = "➖ Neutre"
$recoData[22,0] = "FILTISAC CI (FTSC)"
$recoData[22,1] = 2
$recoData[22,2] = 0
$recoData[22,3] = 14.97
$recoData[22,4] = 7.47
$recoData[22,5] = "🟡 Observer"
$recoData[22,6] = "➖ Neutre"
$recoData[23,0] = "SOLIBRA CI (SLBC)"
$recoData[23,1] = 2
$recoData[23,2] = 0
$recoData[23,3] = 11.44
$recoData[23,4] = 7.48
$recoData[23,5] = "🟡 Observer"
$recoData[23,6] = "➖ Neutre"
$recoData[24,0] = "BERNABE CI (BNBC)"
$recoData[24,1] = 3
$recoData[24,2] = 1
$recoData[24,3] = 10
$recoData[24,4] = 7.32
$recoData[24,5] = "🟢 Achat"
$recoData[24,6] = "✅ Renforcer"
$recoData[25,0] = "CORIS BANK INTERNATIONAL (CBIBF)"
$recoData[25,1] = 1
$recoData[25,2] = 0
$recoData[25,3] = 7.18
$recoData[25,4] = 7.18
$recoData[25,5] = "🟡 Observer"
$recoData[25,6] = "➖ Neutre"
$recoData[26,0] = "SICOR CI (SICC)"
$recoData[26,1] = 1
$recoData[26,2] = 0
$recoData[26,3] = 7.1
$recoData[26,4] = 7.1
$recoData[26,5] = "🟡 Observer"
$recoData[26,6] = "➖ Neutre"
$recoData[27,0] = "SONATEL SN (SNTS)"
$recoData[27,1] = 1
$recoData[27,2] = 0
$recoData[27,3] = 3.59
$recoData[27,4] = 3.59
$recoData[27,5] = "🟡 Observer"
$recoData[27,6] = "➖ Neutre"
$recoData[28,0] = "BANK OF AFRICA SENEGAL (BOAS)"
$recoData[28,1] = 1
$recoData[28,2] = 0
$recoData[28,3] = 3.38
$recoData[28,4] = 3.38
$recoData[28,5] = "🟡 Observer"
$recoData[28,6] = "➖ Neutre"
$recoData[29,0] = "SUCRIVOIRE (SCRC)"
$recoData[29,1] = 1
$recoData[29,2] = 0
$recoData[29,3] = 3.16
$recoData[29,4] = 3.16
$recoData[29,5] = "🟡 Observer"
$recoData[29,6] = "➖ Neutre"
$recoData[30,0] = "SERVAIR ABIDJAN CI (ABJC)"
$recoData[30,1] = 1
$recoData[30,2] = 0
$recoData[30,3] = 2.99
$recoData[30,4] = 2.99
$recoData[30,5] = "🟡 Observer"
$recoData[30,6] = "➖ Neutre"
$recoData[31,0] = "SODE CI (SDCC)"
$recoData[31,1] = 1
$recoData[31,2] = 0
$recoData[31,3] = 1.67
$recoData[31,4] = 1.67
$recoData[31,5] = "🟡 Observer"
$recoData[31,6] = "➖ Neutre"
$recoData[32,0] = "ONATEL BF (ONTBF)"
$recoData[32,1] = 1
$recoData[32,2] = 0
$recoData[32,3] = 1.35
$recoData[32,4] = 1.35
$recoData[32,5] = "🟡 Observer"
$recoData[32,6] = "➖ Neutre"
$recoData[33,0] = "VIVO ENERGY CI (SHEC)"
$recoData[33,1] = 1
$recoData[33,2] = 0
$recoData[33,3] = 0.94
$recoData[33,4] = 0.94
$recoData[33,5] = "🟡 Observer"
$recoData[33,6] = "➖ Neutre"
$recoData[34,0] = "ECOBANK COTE D''IVOIRE (ECOC)"
$recoData[34,1] = 1
$recoData[34,2] = 1
$recoData[34,3] = 0.39
$recoData[34,4] = 3.26
$recoData[34,5] = "🟡 Observer"
$recoData[34,6] = "👀 À surveiller"
$recoData[35,0] = "TRACTAFRIC MOTORS CI (PRSC)"
$recoData[35,1] = 1
$recoData[35,2] = 1
$recoData[35,3] = 0.15
$recoData[35,4] = -3.85
$recoData[35,5] = "🟡 Observer"
$recoData[35,6] = "👀 À surveiller"
$recoData[36,0] = "TOTAL"
$recoData[36,1] = 0
$recoData[36,2] = 4
$recoData[36,3] = 0
$recoData[36,4] = 0
$recoData[36,5] = "🟡 Observer"
$recoData[36,6] = "➖ Neutre"
$recoData[37,0] = "SOGB CI (SOGC)"
$recoData[37,1] = 0
$recoData[37,2] = 1
$recoData[37,3] = -1.22
$recoData[37,4] = -1.22
$recoData[37,5] = "🟡 Observer"
$recoData[37,6] = "➖ Neutre"
$recoData[38,0] = "SAFCA CI (SAFC)"
$recoData[38,1] = 2
$recoData[38,2] = 1
$recoData[38,3] = -1.25
$recoData[38,4] = 3.08
$recoData[38,5] = "🟡 Observer"
$recoData[38,6] = "👀 À surveiller"
$recoData[39,0] = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$recoData[39,1] = 0
$recoData[39,2] = 1
$recoData[39,3] = -1.4
$recoData[39,4] = -1.4
$recoData[39,5] = "🟡 Observer"
$recoData[39,6] = "➖ Neutre"
$recoData[40,0] = "BANK OF AFRICA BN (BOAB)"
$recoData[40,1] = 0
$recoData[40,2] = 1
$recoData[40,3] = -2
$recoData[40,4] = -2
$recoData[40,5] = "🟡 Observer"
$recoData[40,6] = "➖ Neutre"
$recoData[41,0] = "BANK OF AFRICA BF (BOABF)"
$recoData[41,1] = 0
$recoData[41,2] = 1
$recoData[41,3] = -2.14
$recoData[41,4] = -2.14
$recoData[41,5] = "🟡 Observer"
$recoData[41,6] = "➖ Neutre"
$recoData[42,0] = "BICI CI (BICC)"
$recoData[42,1] = 0
$recoData[42,2] = 1
$recoData[42,3] = -2.6
$recoData[42,4] = -2.6
$recoData[42,5] = "🟡 Observer"
$recoData[42,6] = "➖ Neutre"
$recoData[43,0] = "PALM CI (PALC)"
$recoData[43,1] = 0
$recoData[43,2] = 1
$recoData[43,3] = -2.74
$recoData[43,4] = -2.74
$recoData[43,5] = "🟡 Observer"
$recoData[43,6] = "➖ Neutre"
$recoData[44,0] = "BANK OF AFRICA CI (BOAC)"
$recoData[44,1] = 0
$recoData[44,2] = 1
$recoData[44,3] = -3.27
$recoData[44,4] = -3.27
$recoData[44,5] = "🟡 Observer"
$recoData[44,6] = "➖ Neutre"
$recoData[45,0] = "NEI-CEDA CI (NEIC)"
$recoData[45,1] = 0
$recoData[45,2] = 1
$recoData[45,3] = -3.33
$recoData[45,4] = -3.33
$recoData[45,5] = "🟡 Observer"
$recoData[45,6] = "➖ Neutre"
$recoData[46,0] = "ORANGE COTE D'IVOIRE (ORAC)"
$recoData[46,1] = 0
$recoData[46,2] = 1
$recoData[46,3] = -3.33
$recoData[46,4] = -3.33
$recoData[46,5] = "🟡 Observer"
$recoData[46,6] = "➖ Neutre"
$recoData[47,0] = "CFAO MOTORS CI (CFAC)"
$recoData[47,1] = 0
$recoData[47,2] = 2
$recoData[47,3] = -5.15
$recoData[47,4] = -2.21
$recoData[47,5] = "🟡 Observer"
$recoData[47,6] = "➖ Neutre"
$recoData[48,0] = "ECOBANK TRANS. INCORP. TG (ETIT)"
$recoData[48,1] = 0
$recoData[48,2] = 1
$recoData[48,3] = -5.56
$recoData[48,4] = -5.56
$recoData[48,5] = "🟡 Observer"
$recoData[48,6] = "➖ Neutre"
$recoData[49,0] = "SMB CI (SMBC)"
$recoData[49,1] = 0
$recoData[49,2] = 2
$recoData[49,3] = -5.72
$recoData[49,4] = -2.48
$recoData[49,5] = "🟡 Observer"
$recoData[49,6] = "➖ Neutre"
$recoData[50,0] = "BANK OF AFRICA NG (BOAN)"
$recoData[50,1] = 0
$recoData[50,2] = 2
$recoData[50,3] = -7.9
$recoData[50,4] = -1.25
$recoData[50,5] = "🟡 Observer"
$recoData[50,6] = "➖ Neutre"

$wsReco.Range("A2:G52").Value = $recoData

$wsReco.Rows.Item(54).Delete()
$wsReco.Rows.Item(53).Delete()

$wsTopYtd = $wb.Worksheets.Item("Top_YTD")

$topYtdData = New-Object "object[,]" 10,2
$topYtdData[0,0] = "BRVM - SERVICES PUBLICS"
$topYtdData[0,1] = 9683455.41
$topYtdData[1,0] = "SUCRIVOIRE"
$topYtdData[1,1] = 1304285.12
$topYtdData[2,0] = "SAFCA CI"
$topYtdData[2,1] = 376981.22
$topYtdData[3,0] = "CFAO MOTORS CI"
$topYtdData[3,1] = 374627.6
$topYtdData[4,0] = "BRVM - AUTRES SECTEURS"
$topYtdData[4,1] = 334394.66
$topYtdData[5,0] = "NEI-CEDA CI"
$topYtdData[5,1] = 228129.31
$topYtdData[6,0] = "UNIWAX CI"
$topYtdData[6,1] = 220060.37
$topYtdData[7,0] = "SETAO CI"
$topYtdData[7,1] = 196878.32
$topYtdData[8,0] = "AIR LIQUIDE CI"
$topYtdData[8,1] = 166361.44
$topYtdData[9,0] = "BRVM - DISTRIBUTION"
$topYtdData[9,1] = 50640.71

$wsTopYtd.Range("A2:B11").Value = $topYtdData

Write-Output "BRVM data refresh applied."
